$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark two more to-do items as complete (D column = COMPLETE checkbox)
$ws.Range("D8").Value = $true
$ws.Range("D10").Value = $true

# Re-apply the F-column "IF(Dn,1,0)" formulas in the two blocks separated
# by the edited row (33/34 boundary) so Excel groups them as shared
# formulas, matching how the workbook was actually edited/saved.
$ws.Range("F2:F33").Formula = "=IF(D2,1,0)"
$ws.Range("F34:F55").Formula = "=IF(D34,1,0)"

# Restore the viewer's on-screen selection/scroll position
$ws.Range("C7").Select()
